{"js": "// Replace the generated two-digit-by-two-digit multiplication problems\n// (text + answer) throughout the document body with the new set from\n// the commit. Every old value is unique in the document, so a simple\n// exact-text search/replace per pair is unambiguous.\nconst replacements = [\n  [\"96\u00d752=4992\", \"28\u00d761=1708\"],\n  [\"93\u00d722=2046\", \"36\u00d767=2412\"],\n  [\"46\u00d782=3772\", \"69\u00d711=759\"],\n  [\"21\u00d737=777\", \"72\u00d778=5616\"],\n  [\"65\u00d714=910\", \"61\u00d798=5978\"],\n  [\"65\u00d760=3900\", \"90\u00d716=1440\"],\n  [\"55\u00d747=2585\", \"18\u00d734=612\"],\n  [\"91\u00d724=2184\", \"50\u00d746=2300\"],\n  [\"19\u00d718=342\", \"96\u00d789=8544\"],\n  [\"51\u00d773=3723\", \"21\u00d797=2037\"],\n  [\"39\u00d770=2730\", \"41\u00d727=1107\"],\n  [\"32\u00d736=1152\", \"84\u00d776=6384\"],\n  [\"60\u00d752=3120\", \"55\u00d742=2310\"],\n  [\"68\u00d735=2380\", \"16\u00d733=528\"],\n  [\"23\u00d765=1495\", \"90\u00d771=6390\"],\n  [\"37\u00d765=2405\", \"17\u00d757=969\"],\n  [\"44\u00d758=2552\", \"90\u00d795=8550\"],\n  [\"87\u00d739=3393\", \"37\u00d773=2701\"],\n  [\"87\u00d790=7830\", \"19\u00d760=1140\"],\n  [\"84\u00d773=6132\", \"27\u00d737=999\"],\n  [\"42\u00d716=672\", \"98\u00d738=3724\"],\n  [\"86\u00d736=3096\", \"50\u00d715=750\"],\n  [\"19\u00d766=1254\", \"41\u00d730=1230\"],\n  [\"81\u00d728=2268\", \"56\u00d760=3360\"],\n  [\"62\u00d761=3782\", \"55\u00d785=4675\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the generated two-digit-by-two-digit multiplication problems\n# (text + answer) throughout the document body with the new set from\n# the commit. Every old value is unique in the document, so an exact\n# Find/Replace (whole-document scope, MatchCase on, wildcards off) per\n# pair is unambiguous and replaces exactly one occurrence each.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"96\u00d752=4992\", \"28\u00d761=1708\"),\n    @(\"93\u00d722=2046\", \"36\u00d767=2412\"),\n    @(\"46\u00d782=3772\", \"69\u00d711=759\"),\n    @(\"21\u00d737=777\", \"72\u00d778=5616\"),\n    @(\"65\u00d714=910\", \"61\u00d798=5978\"),\n    @(\"65\u00d760=3900\", \"90\u00d716=1440\"),\n    @(\"55\u00d747=2585\", \"18\u00d734=612\"),\n    @(\"91\u00d724=2184\", \"50\u00d746=2300\"),\n    @(\"19\u00d718=342\", \"96\u00d789=8544\"),\n    @(\"51\u00d773=3723\", \"21\u00d797=2037\"),\n    @(\"39\u00d770=2730\", \"41\u00d727=1107\"),\n    @(\"32\u00d736=1152\", \"84\u00d776=6384\"),\n    @(\"60\u00d752=3120\", \"55\u00d742=2310\"),\n    @(\"68\u00d735=2380\", \"16\u00d733=528\"),\n    @(\"23\u00d765=1495\", \"90\u00d771=6390\"),\n    @(\"37\u00d765=2405\", \"17\u00d757=969\"),\n    @(\"44\u00d758=2552\", \"90\u00d795=8550\"),\n    @(\"87\u00d739=3393\", \"37\u00d773=2701\"),\n    @(\"87\u00d790=7830\", \"19\u00d760=1140\"),\n    @(\"84\u00d773=6132\", \"27\u00d737=999\"),\n    @(\"42\u00d716=672\", \"98\u00d738=3724\"),\n    @(\"86\u00d736=3096\", \"50\u00d715=750\"),\n    @(\"19\u00d766=1254\", \"41\u00d730=1230\"),\n    @(\"81\u00d728=2268\", \"56\u00d760=3360\"),\n    @(\"62\u00d761=3782\", \"55\u00d785=4675\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
